$d = $word.ActiveDocument

# 1. Title number fix: "15. Portmanteau" -> "11. Portmanteau"
$d.Content.Find.Execute("15. Portmanteau", $true, $false, $false, $false, $false,
                         $true, 1, $false, "11. Portmanteau", 2)

# 2. Merge split runs / remove spell-check artifacts around "jangry"
$d.Content.Find.Execute("jangry hungry angry", $true, $false, $false, $false, $false,
                         $true, 1, $false, "jangry hungry angry", 2)

# 3. Merge split runs / remove spell-check artifacts around "whism"
$d.Content.Find.Execute("whism antidisestablishmentarianism what", $true, $false, $false, $false, $false,
                         $true, 1, $false, "whism antidisestablishmentarianism what", 2)
